$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the columns that get shuffled between rows
# (D = Fecha, L = Calidad, M = Volumen, N = Precio minimo, O = Precio maximo,
#  P = Precio promedio ponderado, R = Origen, S = Precio $/Kg) for rows 2..7.
$rows = 2..7
$cols = @("D","L","M","N","O","P","R","S")

$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        if ($c -eq "D") {
            $rowData[$c] = $ws.Range("$c$r").Value2()
        } else {
            $rowData[$c] = $ws.Range("$c$r").Value()
        }
    }
    $snapshot[$r] = $rowData
}

# Mapping of new row -> source (old) row, derived from the diff
$mapping = @{
    2 = 6
    3 = 7
    4 = 5
    5 = 2
    6 = 3
    7 = 4
}

foreach ($newRow in $rows) {
    $oldRow = $mapping[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($c in $cols) {
        if ($c -eq "D") {
            $ws.Range("$c$newRow").Value2 = $src[$c]
        } else {
            $ws.Range("$c$newRow").Value = $src[$c]
        }
    }
}
